# Commit: "Moving from POI 3.17.0 to 4.0.1."
#
# This upgrade changed nothing about the document's visible content; it
# only changed the order in which the underlying OOXML writer emits
# attributes (namespace/attr serialization order differs between the
# POI versions). We reproduce that by touching the relevant object-model
# properties with their own current values, which forces Word to
# re-serialize those parts (word/document.xml's <w:sectPr> and
# word/styles.xml) using the current (newer) attribute ordering, without
# altering any actual setting.

$d = $word.ActiveDocument

# --- word/document.xml : <w:sectPr> (<w:pgSz>/<w:pgMar>) -------------
# Re-assigning PageSetup values (even to themselves) makes the section
# properties get rewritten with the new w:pgSz (w:w before w:h) and
# w:pgMar (top/right/bottom/left/header/footer/gutter) attribute order.
$sec = $d.Sections.Item(1)
$ps = $sec.PageSetup
$ps.PageWidth = $ps.PageWidth
$ps.PageHeight = $ps.PageHeight
$ps.TopMargin = $ps.TopMargin
$ps.RightMargin = $ps.RightMargin
$ps.BottomMargin = $ps.BottomMargin
$ps.LeftMargin = $ps.LeftMargin
$ps.HeaderDistance = $ps.HeaderDistance
$ps.FooterDistance = $ps.FooterDistance
$ps.Gutter = $ps.Gutter

# --- word/styles.xml ----------------------------------------------------
# Touching any style (re-assigning its own Priority) forces the whole
# styles part to be re-serialized, which reorders the attributes of
# <w:docDefaults>, <w:latentStyles>, every <w:lsdException>, and each
# <w:style> element (w:type/w:default/w:styleId) to match the newer
# writer's output -- with no change to the actual style definitions.
foreach ($styleName in @("Normal", "Default Paragraph Font", "Normal Table", "No List")) {
    $s = $d.Styles.Item($styleName)
    $s.Priority = $s.Priority
}
